# Burndown Sheet -> version 1.1 update
# Updates the "Burn Down Chart" sheet: new sprint header, re-purposed task
# rows (only the first task row is now populated, rows 2-4 cleared out),
# refreshed hour counts, resized column C, and a new active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burn Down Chart")
$ws.Activate()

# --- Row 1: project banner text (Release 1.0/Sprint 1 -> Release 1.1/Sprint 2)
$ws.Cells.Item(1,1).Value = "Project Title: Minesweeper" + [char]10 + "Release #:1.1" + [char]10 + "Sprint #: 2"

# --- Row 4 (Task 1): re-purposed to "registration/login" work, more hours logged
$ws.Cells.Item(4,2).Value = 1
$ws.Cells.Item(4,3).Value = "Create registration page and login page"
$ws.Cells.Item(4,4).Value = "Create Registration page"
$ws.Cells.Item(4,5).Value = "Team"
$ws.Cells.Item(4,6).Value = 5
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = 0
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(4,10).Value = 0
$ws.Cells.Item(4,11).Value = 0
$ws.Cells.Item(4,12).Value = 2
$ws.Cells.Item(4,13).Value = 2
$ws.Rows.Item(4).RowHeight = 30

# --- Row 5 (was Task 2): cleared except for the task-description cell
$ws.Cells.Item(5,2).Value = ""
$ws.Cells.Item(5,3).Value = ""
$ws.Cells.Item(5,4).Value = "Create login page"
$ws.Cells.Item(5,5).Value = ""
$ws.Cells.Item(5,6).Value = ""
$ws.Cells.Item(5,7).Value = ""
$ws.Cells.Item(5,8).Value = ""
$ws.Cells.Item(5,9).Value = ""
$ws.Cells.Item(5,10).Value = ""
$ws.Cells.Item(5,11).Value = ""
$ws.Cells.Item(5,12).Value = ""
$ws.Cells.Item(5,13).Value = ""
$ws.Rows.Item(5).AutoFit()

# --- Row 6 (was Task 3): cleared except for the task-description cell
$ws.Cells.Item(6,2).Value = ""
$ws.Cells.Item(6,3).Value = ""
$ws.Cells.Item(6,4).Value = "Create Database"
$ws.Cells.Item(6,5).Value = ""
$ws.Cells.Item(6,6).Value = ""
$ws.Cells.Item(6,7).Value = ""
$ws.Cells.Item(6,8).Value = ""
$ws.Cells.Item(6,9).Value = ""
$ws.Cells.Item(6,10).Value = ""
$ws.Cells.Item(6,11).Value = ""
$ws.Cells.Item(6,12).Value = ""
$ws.Cells.Item(6,13).Value = ""
$ws.Rows.Item(6).AutoFit()

# --- Row 7 (was Task 4): cleared except for the task-description cell
$ws.Cells.Item(7,2).Value = ""
$ws.Cells.Item(7,3).Value = ""
$ws.Cells.Item(7,4).Value = "refactoring"
$ws.Cells.Item(7,5).Value = ""
$ws.Cells.Item(7,6).Value = ""
$ws.Cells.Item(7,7).Value = ""
$ws.Cells.Item(7,8).Value = ""
$ws.Cells.Item(7,9).Value = ""
$ws.Cells.Item(7,10).Value = ""
$ws.Cells.Item(7,11).Value = ""
$ws.Cells.Item(7,12).Value = ""
$ws.Cells.Item(7,13).Value = ""
$ws.Rows.Item(7).AutoFit()

# --- Column C: narrower, best-fit width (~18.86 chars; closest attainable
#     snap point in this engine's pixel-quantised column-width storage)
$ws.Columns.Item(3).ColumnWidth = 18
$ws.Columns.Item(3).BestFit = $true

# --- Recalculate so the "Actual Remaining" row (9) and chart feeders refresh
$wb.Application.Calculate()

# --- Selection moves to M8 (matches the saved view in the edited workbook)
$ws.Range("M8").Select()
